$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Footer "datetimeFigureOut" field: 11/15/2010 -> 5/27/2011
#    (slide master + every custom layout)
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    foreach ($sh in $container.Shapes) {
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "5/27/2011"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout
}

# ---------------------------------------------------------------------------
# 2. Slide 1 - "TextBox 15": "is descriptor of" -> 3 runs "is " / "descriptor " / "of"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
foreach ($sh in $s1.Shapes) {
    if ($sh.Name -eq "TextBox 15") {
        $tr = $sh.TextFrame.TextRange
        $tr.Characters(1, 3).Text = "is "
        $tr.Characters(4, 11).Text = "descriptor "
        $tr.Characters(15, 2).Text = "of"
    }
}

# ---------------------------------------------------------------------------
# 3. Slide 4 - "TextBox 5": drop stray empty run, merge "data "+"format specification"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
foreach ($sh in $s4.Shapes) {
    if ($sh.Name -eq "TextBox 5") {
        $tr = $sh.TextFrame.TextRange
        $tr.Characters(1, 9).Text = "numeric "
        $tr.Characters(10, 25).Text = "data format specification"
    }
    if ($sh.Name -eq "TextBox 19") {
        $tr = $sh.TextFrame.TextRange
        $tr.Characters(1, 9).Text = "textual "
    }
}

# ---------------------------------------------------------------------------
# 4. Slide 5 - "TextBox 31": merge "data "+"item" -> "data item"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
foreach ($sh in $s5.Shapes) {
    if ($sh.Name -eq "TextBox 31") {
        $tr = $sh.TextFrame.TextRange
        $tr.Characters(1, 9).Text = "data item"
    }
}
